$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.420.37'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.020.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.83%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.662'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.91%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '45.34'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.81'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.04%  '

$ws.Range("E10").Value = '  +1.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0717'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.18%  '

$ws.Range("E12").Value = '  +0.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.71'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.316.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.88%  '

$ws.Range("E15").Value = '  +1.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.025.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.532.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.94%  '

$ws.Range("E20").Value = '  -0.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '236.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.07%  '

$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.97%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.63'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.29%  '

$ws.Range("E29").Value = '  -9.17%  '

$ws.Range("E30").Value = '  -4.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +56.83%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.54%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0593'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.75%  '

$ws.Range("E34").Value = '  +0.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.86'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.71%  '

$ws.Range("E36").Value = '  -4.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0815'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +12.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.841'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.02%  '

$ws.Range("E40").Value = '  -8.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0217'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.51%  '

$ws.Range("E43").Value = '  +2.19%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.82%  '

$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.74'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +13.82%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.320.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.74%  '

$ws.Range("E47").Value = '  +0.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.226.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.62%  '
